$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# New Madigan bike hours: update Riders and Average columns for Friday (row 2) and Saturday (row 3)
$ws.Range("C2").Value = 281
$ws.Range("D2").Value = 281

$ws.Range("C3").Value = 127
$ws.Range("D3").Value = 127
